# Mark DRAG-related requirements as tested / Complete:
#   - Leading Edge Flap Deflection (SYS-LFD-xxx) -> rows 116-122
#   - Tailing Flap Deflection      (SYS-TFD-xxx) -> rows 229-232
#   - Gear                         (SYS-GDR-xxx) -> rows 242-246
#   - Speedbrake                   (SYS-SDE-xxx) -> rows 247-251

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats constant
$xlPasteFormats = -4122

# --- Leading Edge Flap Deflection: E116:E122 ---------------------------
# Reuse the "Complete" formatting (style) already present on E111 so the
# resulting cell style exactly matches the other "Complete" rows.
$fmtSrcLfd = $ws.Range("E111")
$fmtSrcLfd.Copy()
$ws.Range("E116:E122").PasteSpecial($xlPasteFormats)
$ws.Range("E116:E122").Value = "Complete"

# --- Tailing Flap Deflection: E229:E232 --------------------------------
# --- Gear: E242:E246 -----------------------------------------------------
# --- Speedbrake: E247:E251 -----------------------------------------------
# These all share the other "Complete" style, already present on E139.
$fmtSrcOther = $ws.Range("E139")
$fmtSrcOther.Copy()
$ws.Range("E229:E232").PasteSpecial($xlPasteFormats)
$ws.Range("E229:E232").Value = "Complete"

$fmtSrcOther.Copy()
$ws.Range("E242:E251").PasteSpecial($xlPasteFormats)
$ws.Range("E242:E251").Value = "Complete"

# --- Update the saved view / selection ----------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 100
$win.ScrollColumn = 4
$ws.Range("D114").Select()

$excel.CutCopyMode = $false
